# Version 3.0 - added changing column names
# Function to change the column headers (Column 1 / New Column 2) and
# populate the second column with sequential letters (A, B, C, ...) for
# every data row, applied to every worksheet in the workbook.

function Set-ColumnHeaders {
    param(
        $Worksheet
    )

    # Update the header row.
    $Worksheet.Range("A1").Value = "Column 1"
    $Worksheet.Range("B1").Value = "New Column 2"

    # Work out how many data rows exist (based on column A).
    $lastRow = $Worksheet.Cells.Item(1, 1).End(-4121).Row

    # Fill column B with sequential letters A, B, C, ... for each data row.
    for ($row = 2; $row -le $lastRow; $row++) {
        $letter = [char](65 + ($row - 2))
        $Worksheet.Cells.Item($row, 2).Value = [string]$letter
    }
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    Set-ColumnHeaders $ws
}
